$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    paragraph ("Play Down the Rails Free - Review of Pragmatic Play's Slot
#    Game"). The new paragraph has: an empty leading run, a bold run with
#    "Meta description", and a plain run with the rest of the sentence.
# ---------------------------------------------------------------------------

$titlePara = $d.Paragraphs.Item(1)
$titleRange = $titlePara.Range
$titleRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaRange = $metaPara.Range

$metaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
<pkg:xmlData>
<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships>
</pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r/>
<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>
<w:r><w:t>: Read our review of Down the Rails, the London subway-themed slot game from Pragmatic Play. Play Down the Rails for free and enjoy bonus games and random features.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$metaRange.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2) Near the end of the document, remove the duplicate bold paragraph
#    ("Play Down the Rails Free - Review of Pragmatic Play's Slot Game")
#    and replace the text of the final (italic) paragraph with the new
#    DALL-E image prompt, keeping its italic formatting.
# ---------------------------------------------------------------------------

$total = $d.Paragraphs.Count
$dupTitlePara = $null
for ($i = $total; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Play Down the Rails Free - Review of Pragmatic Play's Slot Game`r") {
        $dupTitlePara = $p
        break
    }
}
$dupTitlePara.Range.Delete()

$total2 = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($total2)
$lastRange = $lastPara.Range
# Exclude the trailing paragraph-mark character from the replaced text.
$textRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$textRange.Text = "Prompt for DALLE: Create a feature image for ""Down the Rails"" that captures the game's theme of the London subway and features a happy Maya warrior with glasses. The image should be in a cartoon style and include elements from the game such as the London subway, iconic characters like Shakespeare or Winston Churchill, and bonus features like the End of the Line Bonus game. The Maya warrior should be prominently displayed, perhaps riding the subway or standing in front of Buckingham Palace. The image should be eye-catching and colorful, with a sense of fun and excitement to match the game."
